$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 0.9811601196580918
$ws.Range("BP1").Value = 0.93049165827590952
$ws.Range("A2").Value = 0.77580724980752813
$ws.Range("D2").Value = 0.76223124490013561
$ws.Range("BP2").Value = 0.78917830149659207
$ws.Range("B3").Value = 0.63907073204947362
$ws.Range("V3").Value = 0.91227132869627114
$ws.Range("C5").Value = 0.97452961514690473
$ws.Range("F5").Value = 0.59758593918418501
$ws.Range("D6").Value = 0.87320026544191176
$ws.Range("E7").Value = 0.85526944704841479
$ws.Range("F7").Value = 0.94892009155707435
$ws.Range("G8").Value = 0.85765176500167839
$ws.Range("G9").Value = 0.95238020131827383
$ws.Range("K9").Value = 0.96832258681677652
$ws.Range("H10").Value = 0.8042736446165426
$ws.Range("I10").Value = 0.71317889394818601
$ws.Range("F11").Value = 0.82884980258872742
$ws.Range("J11").Value = 0.86892764092634467
$ws.Range("V11").Value = 0.9296230795448075
$ws.Range("L13").Value = 0.83537939909408776
$ws.Range("B14").Value = 0.76397508179366869
$ws.Range("J14").Value = 0.92774004978385527
$ws.Range("L14").Value = 0.77734087085826553
$ws.Range("AA14").Value = 0.93073947845747451
$ws.Range("F15").Value = 0.99960469088814308
$ws.Range("M15").Value = 0.97950787499907921
$ws.Range("N15").Value = 0.75727075362067642
$ws.Range("Q15").Value = 0.89264036866478902
$ws.Range("BI16").Value = 0.8625726215915841
$ws.Range("P17").Value = 0.91837301366792334
$ws.Range("S17").Value = 0.92630010687754738
$ws.Range("P18").Value = 0.84745787227355418
$ws.Range("Q18").Value = 0.90511986346440443
$ws.Range("S18").Value = 0.82862081351744687
$ws.Range("T18").Value = 0.52106677290557768
$ws.Range("BD18").Value = 0.90361888390836598
$ws.Range("T19").Value = 0.51949561337232275
$ws.Range("U19").Value = 0.67065730282689007
$ws.Range("V20").Value = 0.75149412252086489
$ws.Range("V21").Value = 0.80298159120256918
$ws.Range("W21").Value = 0.90441495442881248
$ws.Range("L23").Value = 0.79832963905269672
$ws.Range("V23").Value = 0.91824266758032158
$ws.Range("V24").Value = 0.98024396303251082
$ws.Range("W25").Value = 0.6761871392911758
$ws.Range("Z25").Value = 0.69097955528154642
$ws.Range("AA25").Value = 0.59728114105203201
$ws.Range("E26").Value = 0.95430484005366334
$ws.Range("Z27").Value = 0.97481533237734974
$ws.Range("AA28").Value = 0.82062999167060213
$ws.Range("BG28").Value = 0.75556373734932558
$ws.Range("AE29").Value = 0.96851713497441549
$ws.Range("AC30").Value = 0.76355538005975365
$ws.Range("BC30").Value = 0.88928286861649375
$ws.Range("BD30").Value = 0.75806744760135714
$ws.Range("M31").Value = 0.98306459256651046
$ws.Range("AS31").Value = 0.92500958106403852
$ws.Range("AE32").Value = 0.90816037399848604
$ws.Range("AG32").Value = 0.85177307478118136
$ws.Range("AI33").Value = 0.84894917488629118
$ws.Range("O34").Value = 0.8623830003939521
$ws.Range("AJ34").Value = 0.7960856010860895
$ws.Range("AH35").Value = 0.61491332718429792
$ws.Range("AJ35").Value = 0.87777894102342224
$ws.Range("AK35").Value = 0.66081358933869438
$ws.Range("X37").Value = 0.96607526541938493
$ws.Range("AJ37").Value = 0.58922971266596225
$ws.Range("AL37").Value = 0.79118002822396827
$ws.Range("AM37").Value = 0.9626672622520267
$ws.Range("J38").Value = 0.80150218151285557
$ws.Range("AJ38").Value = 0.77725935457090234
$ws.Range("AS38").Value = 0.91051016801473628
$ws.Range("BE38").Value = 0.92768590270754614
$ws.Range("AN39").Value = 0.95050514651407658
$ws.Range("BP39").Value = 0.89078906636749444
$ws.Range("AP40").Value = 0.7821657685468979
$ws.Range("AN41").Value = 0.84595425816206848
$ws.Range("AP41").Value = 0.90081215763731359
$ws.Range("AQ42").Value = 0.95982139319267734
$ws.Range("AS42").Value = 0.66437897386374689
$ws.Range("AB43").Value = 0.71390757115003911
$ws.Range("AL43").Value = 0.88125151277774783
$ws.Range("AO44").Value = 0.72605111947636625
$ws.Range("AQ44").Value = 0.97655526401450676
$ws.Range("AS44").Value = 0.94806507951923535
$ws.Range("BN44").Value = 0.82641746595946419
$ws.Range("AT45").Value = 0.99133393506134926
$ws.Range("AH46").Value = 0.91645491246673516
$ws.Range("AS47").Value = 0.61775422170727912
$ws.Range("AV47").Value = 0.89275396889809255
$ws.Range("AW47").Value = 0.7631058399246915
$ws.Range("AT48").Value = 0.74454229682384176
$ws.Range("AW48").Value = 0.87284602406716205
$ws.Range("AX48").Value = 0.92424000000394724
$ws.Range("A49").Value = 0.85235611279686663
$ws.Range("H49").Value = 0.80169639595963949
$ws.Range("AY49").Value = 0.73785634822021762
$ws.Range("AF50").Value = 0.57377442424144398
$ws.Range("AY50").Value = 0.83975031564942193
$ws.Range("BO51").Value = 0.98027398624402529
$ws.Range("Y52").Value = 0.72276186346140148
$ws.Range("AC52").Value = 0.98827977063021244
$ws.Range("AX52").Value = 0.94038045403084647
$ws.Range("AY53").Value = 0.8697325084862827
$ws.Range("BA54").Value = 0.97042443651324728
$ws.Range("AI55").Value = 0.8346838618446133
$ws.Range("BA55").Value = 0.67559918230020821
$ws.Range("BB55").Value = 0.81999040926664102
$ws.Range("D56").Value = 0.77552984195835928
$ws.Range("AG56").Value = 0.7139357499703376
$ws.Range("AO56").Value = 0.74929719690449814
$ws.Range("BB56").Value = 0.87841657614734925
$ws.Range("BC57").Value = 0.96183672226651762
$ws.Range("BE58").Value = 0.85022961267051844
$ws.Range("BG58").Value = 0.85708498216777751
$ws.Range("BH58").Value = 0.94334044540144601
$ws.Range("BH59").Value = 0.84143432761112713
$ws.Range("BI59").Value = 0.97003156972527826
$ws.Range("BI60").Value = 0.6182967660375468
$ws.Range("BJ60").Value = 0.95079386513254627
$ws.Range("T62").Value = 0.75480858314986876
$ws.Range("BL62").Value = 0.87147423912788757
$ws.Range("BI63").Value = 0.88694171198803429
$ws.Range("BM63").Value = 0.92624902612836202
$ws.Range("BP63").Value = 0.99475517384000223
$ws.Range("BL65").Value = 0.80254336285545247
$ws.Range("BN65").Value = 0.83660563589007486
$ws.Range("BO65").Value = 0.93671328112959462
$ws.Range("BL66").Value = 0.96935568923703186
$ws.Range("BP66").Value = 0.79324343045082757
$ws.Range("A67").Value = 0.81718408808870158
$ws.Range("T67").Value = 0.81336993139075353
$ws.Range("BN67").Value = 0.80829069354328609
$ws.Range("X68").Value = 0.83508116517173048
$ws.Range("BO68").Value = 0.72516788449809688

Write-Output "Applied 136 cell updates"
